# Update cryptos list data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values look like plain numbers (e.g. "0.9998", "244.64").
# Force those specific cells to Text format *before* assigning the string
# so Excel doesn't silently reinterpret/round them as numeric values.
# (Values containing two dots, like "30.380.26", are never auto-parsed as
# numbers, so they don't need this treatment.)
$numericLookingPriceCells = @(
  "D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16",
  "D19","D20","D22","D23","D24","D25","D26","D28","D29","D30","D31",
  "D32","D33","D34","D36","D37","D38","D39","D40","D41","D42","D44",
  "D45","D46","D47","D48","D49","D50","D51"
)
foreach ($addr in $numericLookingPriceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.380.26"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.880.68"
$ws.Range("E3").Value = "  +0.87%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "244.64"
$ws.Range("E5").Value = "  +4.42%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.0000"

# Row 7 - XRP
$ws.Range("D7").Value = "0.4777"
$ws.Range("E7").Value = "  +1.98%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2883"
$ws.Range("E8").Value = "  +1.33%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06525"
$ws.Range("E9").Value = "  -0.24%  "

# Row 10 - Solana
$ws.Range("D10").Value = "21.38"
$ws.Range("E10").Value = "  +0.20%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07765"
$ws.Range("E11").Value = "  +0.21%  "

# Row 12 - was Polygon, now WrappedEther
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.891.24"
$ws.Range("E12").Value = "  +1.28%  "

# Row 13 - Litecoin
$ws.Range("D13").Value = "96.64"
$ws.Range("E13").Value = "  +1.04%  "

# Row 14 - was WrappedEther, now Polygon
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.7365"
$ws.Range("E14").Value = "  +7.01%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "5.136"
$ws.Range("E15").Value = "  +1.15%  "

# Row 16 - BitcoinCash
$ws.Range("D16").Value = "275.48"
$ws.Range("E16").Value = "  +3.81%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "30.360.91"
$ws.Range("E17").Value = "  +0.63%  "

# Row 18 - Avalanche
$ws.Range("E18").Value = "  -1.54%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.000007553"
$ws.Range("E19").Value = "  -1.86%  "

# Row 20 - Dai
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.03%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.125.91"
$ws.Range("E21").Value = "  +0.00%  "

# Row 22 - BinanceUSD
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.01%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.235"
$ws.Range("E23").Value = "  +0.16%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "6.174"
$ws.Range("E24").Value = "  +0.53%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "9.263"
$ws.Range("E25").Value = "  -2.04%  "

# Row 26 - Monero
$ws.Range("D26").Value = "163.66"
$ws.Range("E26").Value = "  -1.39%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +1.76%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "1.962"
$ws.Range("E28").Value = "  +1.79%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "1.371"
$ws.Range("E29").Value = "  +0.28%  "

# Row 30 - Stellar
$ws.Range("D30").Value = "0.09982"
$ws.Range("E30").Value = "  +0.81%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "1.512"
$ws.Range("E31").Value = "  +3.57%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "4.317"
$ws.Range("E32").Value = "  -0.44%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "4.087"
$ws.Range("E33").Value = "  +1.48%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.04741"
$ws.Range("E34").Value = "  +0.38%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  -0.17%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.6965"
$ws.Range("E36").Value = "  +0.10%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01860"
$ws.Range("E38").Value = "  +0.26%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "2.748"
$ws.Range("E39").Value = "  -0.73%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "6.276"
$ws.Range("E40").Value = "  -0.47%  "

# Row 41 - was TheSandbox, now Aave
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "69.51"
$ws.Range("E41").Value = "  -3.52%  "

# Row 42 - was Aave, now TheSandbox
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.4173"
$ws.Range("E42").Value = "  +1.18%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  -0.73%  "

# Row 44 - TrustWalletToken
$ws.Range("D44").Value = "0.8415"
$ws.Range("E44").Value = "  +1.03%  "

# Row 45 - PaxDollar
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  -0.03%  "

# Row 46 - Quant
$ws.Range("D46").Value = "101.89"
$ws.Range("E46").Value = "  -0.81%  "

# Row 47 - was EnergySwap, now Aptos
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.096"
$ws.Range("E47").Value = "  +0.12%  "

# Row 48 - was Aptos, now EnergySwap
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.204"
$ws.Range("E48").Value = "  +1.38%  "

# Row 49 - Elrond
$ws.Range("D49").Value = "35.24"
$ws.Range("E49").Value = "  +2.08%  "

# Row 50 - Maker
$ws.Range("D50").Value = "912.12"
$ws.Range("E50").Value = "  -5.45%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.05591"
$ws.Range("E51").Value = "  -0.86%  "
